# Rename the "pathology" column header to "pathologyID" (used consistently
# across all metadata files per the commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "pathologyID"

# Widen column D slightly to fit the longer header text.
$ws.Columns.Item(4).ColumnWidth = 12.13

# Default column width grew marginally too (sheet-wide formatting nudge).
$ws.StandardWidth = 11.83984375

# Move the active cell/selection from C5 to A4.
$ws.Range("A4").Select() | Out-Null
